$wb = $excel.ActiveWorkbook

$values = @{
    3  = 1395
    4  = 6865
    5  = 388
    6  = 203
    7  = 3928
    8  = 49
    9  = 30
    10 = 50
    11 = 871
    12 = 251
    13 = 5424
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $values.Keys) {
        $ws.Cells.Item($row, 6).Value = $values[$row]
    }
}
